$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append the new day's row of data (row 31) --------------------------
# Pull number formatting (date / 2-decimal) down from row 30 first, mirroring
# what a user does when continuing a table by filling the row below.
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Range("N30:R30").Copy()
$ws.Range("N31:R31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A31").Value = 44123
$ws.Range("B31").Value = 30
$ws.Range("D31").Value = 166
$ws.Range("E31").Value = 113
$ws.Range("F31").Value = 1533
$ws.Range("G31").Value = 893
$ws.Range("H31").Value = 93
$ws.Range("I31").Value = 20
$ws.Range("J31").Value = 14
$ws.Range("K31").Value = 213
$ws.Range("L31").Value = 86

$ws.Range("N31").Formula = "=100*E31/D31"
$ws.Range("O31").Formula = "=100*G31/F31"
$ws.Range("P31").Formula = "=100*H31/D31"
$ws.Range("Q31").Formula = "=100*J31/I31"
$ws.Range("R31").Formula = "=100*L31/K31"

$excel.CalculateFullRebuild()

# --- Extend the two charts so their series cover the new row ------------
$chart1 = $ws.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = "=SERIES(Sheet1!`$N`$1,Sheet1!`$B`$2:`$B`$31,Sheet1!`$N`$2:`$N`$31,1)"
$chart1.SeriesCollection().Item(2).Formula = "=SERIES(Sheet1!`$O`$1,Sheet1!`$B`$2:`$B`$31,Sheet1!`$O`$2:`$O`$31,2)"
$chart1.SeriesCollection().Item(3).Formula = "=SERIES(Sheet1!`$P`$1,Sheet1!`$B`$2:`$B`$31,Sheet1!`$P`$2:`$P`$31,3)"
$chart1.SeriesCollection().Item(4).Formula = "=SERIES(Sheet1!`$Q`$1,Sheet1!`$B`$2:`$B`$31,Sheet1!`$Q`$2:`$Q`$31,4)"
$chart1.SeriesCollection().Item(5).Formula = "=SERIES(Sheet1!`$R`$1,Sheet1!`$B`$2:`$B`$31,Sheet1!`$R`$2:`$R`$31,5)"
$chart1.Refresh()

$chart2 = $ws.ChartObjects().Item(2).Chart
$chart2.SeriesCollection().Item(1).Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$31,Sheet1!`$N`$2:`$N`$31,1)"
$chart2.SeriesCollection().Item(2).Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$31,Sheet1!`$O`$2:`$O`$31,2)"
$chart2.SeriesCollection().Item(3).Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$31,Sheet1!`$P`$2:`$P`$31,3)"
$chart2.SeriesCollection().Item(4).Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$31,Sheet1!`$Q`$2:`$Q`$31,4)"
$chart2.SeriesCollection().Item(5).Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$31,Sheet1!`$R`$2:`$R`$31,5)"
$chart2.Refresh()

# --- Restore the window/selection state ----------------------------------
$ws.Range("O41").Select()
